$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FPVA25")
Write-Host $ws.Name
